# The commit swaps the "Integral" theme colours and the "Office Theme"
# colours between the deck's two theme parts: the slide master's theme
# (which carried the "Integral" palette) picks up the plain "Office"
# palette that used to live in the (otherwise unused) notes-master theme,
# and vice versa. The net effect reachable from the slide/master's own
# theme is: every theme colour slot on ActivePresentation.SlideMaster
# flips from the green "Integral" swatches to the default blue "Office"
# swatches.
#
# PowerPoint's ThemeColorScheme.Colors(i).RGB takes a standard VBA RGB
# long (R + G*256 + B*65536), so convert each target hex triplet before
# assigning it.

function ConvertTo-RgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = ConvertTo-RgbLong $officeColors[$i - 1]
}
